$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page setup used by the other quarterly sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q1" sheet with the fund holding data, matching
#    the layout used by the other quarterly sheets (e.g. "2021-Q4").
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row: A2 numeric index, H2 numeric rank; the rest stored as text.
$newSheet.Range("A2").Value = 0
# Force these to stay text (the source strings look numeric, e.g. "007280"
# and "1.35", and must round-trip as inline strings, not numbers).
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "007280"
$newSheet.Range("C2").Value = "上投摩根日本精选股票（QDII）"
$newSheet.Range("D2").Value = "1.35"
$newSheet.Range("E2").Value = "88.71"
$newSheet.Range("F2").Value = "5.92"
$newSheet.Range("G2").Value = "0.0799"
$newSheet.Range("H2").Value = 1

# Drop the "Text" number-format styling we used to pin these as strings so
# the cells fall back to the plain/default style (matching the sibling
# quarterly sheets), then re-apply the bold/bordered header style.
$newSheet.Range("B2:G2").ClearFormats()

# Re-apply the bold/bordered header style (and the A2 index-cell style) by
# copying it over from the "2021-Q4" sheet, which uses the same layout.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row for "2022-Q1" at
#    the top of the data and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$sumSheet = $wb.Worksheets.Item("总计")
$sumSheet.Rows(2).Insert()

# Give the newly inserted row's index cell (A2) the same bold/bordered
# style as the rest of the index column, then clear the formatting that
# the insert leaked into B2:D2 from the row above.
$sumSheet.Range("A3").Copy()
$sumSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$sumSheet.Range("B2:D2").ClearFormats()

$sumSheet.Range("A2").Value = 0
$sumSheet.Range("B2").Value = "2022-Q1"
$sumSheet.Range("C2").Value = 1
$sumSheet.Range("D2").Value = 0.08

$sumSheet.Range("A3").Value = 1
$sumSheet.Range("B3").Value = "2021-Q4"
$sumSheet.Range("C3").Value = 1
$sumSheet.Range("D3").Value = 0.08

$sumSheet.Range("A4").Value = 2
$sumSheet.Range("B4").Value = "2021-Q3"
$sumSheet.Range("C4").Value = 1
$sumSheet.Range("D4").Value = 0.07000000000000001

$sumSheet.Range("A5").Value = 3
$sumSheet.Range("B5").Value = "2021-Q2"
$sumSheet.Range("C5").Value = 1
$sumSheet.Range("D5").Value = 0.06

$sumSheet.Range("A6").Value = 4
$sumSheet.Range("B6").Value = "2021-Q1"
$sumSheet.Range("C6").Value = 1
$sumSheet.Range("D6").Value = 0.06

$sumSheet.Range("A7").Value = 5
$sumSheet.Range("B7").Value = "2020-Q4"
$sumSheet.Range("C7").Value = 1
$sumSheet.Range("D7").Value = 0.05

# ---------------------------------------------------------------------------
# 4. Restore the originally active sheet/tab ("2020-Q4") since adding a new
#    worksheet shifts Excel's focus onto it.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
